$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "Date" before "Screenshot Path", shifting
# "Screenshot Path" to the new column P ---
$ws.Range("O1").Value = "Date"
$ws.Range("P1").Value = "Screenshot Path"

# Give the new header cell P1 the same formatting (bold, centered, bordered)
# as the rest of the header row by copying format from its neighbor O1.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("P1").Value = "Screenshot Path"

# --- Data row (row 2): replace sample data with placeholder tokens, and
# move the screenshot path value to the new column P ---
$ws.Range("A2").Value = "A2"
$ws.Range("B2").Value = "A3"
$ws.Range("C2").Value = "A4"
$ws.Range("D2").Value = "A5"
$ws.Range("E2").Value = "A6"
$ws.Range("F2").Value = "A7"
$ws.Range("G2").Value = "A8"
$ws.Range("H2").Value = "A9"
$ws.Range("I2").Value = "A10"
$ws.Range("J2").Value = "A11"
$ws.Range("K2").Value = "A12"
$ws.Range("L2").Value = "A13"
$ws.Range("M2").Value = "A14"
$ws.Range("N2").Value = "A15"
$ws.Range("O2").Value = "A16"
$ws.Range("P2").Value = "C:/Users/seema/Pictures/Screenshots/Screenshots(2).png"

# Match the saved selection state (whole second row selected, active cell A2)
$ws.Range("A2:XFD2").Select() | Out-Null
